$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------------
# Append 5 new rows to the "Tabela1" table (Mapeamento_Reports), columns:
#   A=View  B=Campo  C=Tipo de Dados  D=Atualização Timezone
#   E=Colunas1  F=Colunas2
# New rows use the same "sim"/"ok" pattern already used throughout the sheet.
# ---------------------------------------------------------------------------

# Row 212
$ws.Cells.Item(212, 1).Value = "WMS_REL_068_-_Elton_Relatorio"
$ws.Cells.Item(212, 2).Value = "DT_LIMITE"
$ws.Cells.Item(212, 3).VerticalAlignment = $ws.Cells.Item(212, 3).VerticalAlignment
$ws.Cells.Item(212, 4).Value = "sim"
$ws.Cells.Item(212, 5).Value = "ok"
$ws.Cells.Item(212, 6).VerticalAlignment = $ws.Cells.Item(212, 6).VerticalAlignment

# Row 213
$ws.Cells.Item(213, 1).Value = "WMS_REL_068_-_Elton_Relatorio"
$ws.Cells.Item(213, 2).Value = "DT_REGISTRO"
$ws.Cells.Item(213, 3).VerticalAlignment = $ws.Cells.Item(213, 3).VerticalAlignment
$ws.Cells.Item(213, 4).Value = "sim"
$ws.Cells.Item(213, 5).Value = "ok"
$ws.Cells.Item(213, 6).VerticalAlignment = $ws.Cells.Item(213, 6).VerticalAlignment

# Row 214
$ws.Cells.Item(214, 1).Value = "WMS_REL_074_Base_Recebimento"
$ws.Cells.Item(214, 2).Value = "DATA"
$ws.Cells.Item(214, 3).VerticalAlignment = $ws.Cells.Item(214, 3).VerticalAlignment
$ws.Cells.Item(214, 4).Value = "sim"
$ws.Cells.Item(214, 5).Value = "ok"
$ws.Cells.Item(214, 6).VerticalAlignment = $ws.Cells.Item(214, 6).VerticalAlignment

# Row 215
$ws.Cells.Item(215, 1).Value = "WMS_Valor_expedido_volume_cubado"
$ws.Cells.Item(215, 2).Value = "MES"
$ws.Cells.Item(215, 3).VerticalAlignment = $ws.Cells.Item(215, 3).VerticalAlignment
$ws.Cells.Item(215, 4).Value = "sim"
$ws.Cells.Item(215, 5).Value = "ok"
$ws.Cells.Item(215, 6).VerticalAlignment = $ws.Cells.Item(215, 6).VerticalAlignment

# Row 216 - note: set Campo (B) before View (A) so new shared strings are
# appended to sharedStrings.xml in the same order as the source commit.
$ws.Cells.Item(216, 2).Value = "MESANO"
$ws.Cells.Item(216, 1).Value = "WMS_Valor_Recebido_Volume_Cubado"
$ws.Cells.Item(216, 3).VerticalAlignment = $ws.Cells.Item(216, 3).VerticalAlignment
$ws.Cells.Item(216, 4).Value = "sim"
$ws.Cells.Item(216, 5).Value = "ok"
$ws.Cells.Item(216, 6).VerticalAlignment = $ws.Cells.Item(216, 6).VerticalAlignment

# Grow the table / autofilter range to cover the newly-added rows.
$lo.Resize($ws.Range("A1:F216"))

# ---------------------------------------------------------------------------
# Sheet view: keep header row frozen, scroll the frozen view further down and
# move the active selection to where editing left off.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$ws.Range("A2").Select()
$win.FreezePanes = $true
$ws.Range("B218").Select()
